$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update column F ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 164
$wsExhibit.Range("F4").Value = 8077
$wsExhibit.Range("F14").Value = 73
$wsExhibit.Range("F17").Value = 5992
$wsExhibit.Range("F20").Value = 2106
$wsExhibit.Range("F21").Value = 62
$wsExhibit.Range("F22").Value = 102
$wsExhibit.Range("F24").Value = 425

# Sheet "全部类型" (all types) — same underlying data, update column F values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 164
$wsAll.Range("F4").Value = 8077
$wsAll.Range("F15").Value = 73
$wsAll.Range("F19").Value = 5992
$wsAll.Range("F23").Value = 2106
$wsAll.Range("F24").Value = 62
$wsAll.Range("F25").Value = 102
$wsAll.Range("F27").Value = 425
